$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-12 Friday", 2)
$d.Content.Find.Execute("899÷2=449, 1", $true, $false, $false, $false, $false, $true, 1, $false, "551÷4=137, 3", 2)
$d.Content.Find.Execute("498÷7=71, 1", $true, $false, $false, $false, $false, $true, 1, $false, "281÷7=40, 1", 2)
$d.Content.Find.Execute("353÷5=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "199÷8=24, 7", 2)
$d.Content.Find.Execute("563÷8=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "147÷7=21, 0", 2)
$d.Content.Find.Execute("176÷9=19, 5", $true, $false, $false, $false, $false, $true, 1, $false, "919÷2=459, 1", 2)
$d.Content.Find.Execute("955÷6=159, 1", $true, $false, $false, $false, $false, $true, 1, $false, "992÷3=330, 2", 2)
$d.Content.Find.Execute("889÷9=98, 7", $true, $false, $false, $false, $false, $true, 1, $false, "270÷5=54, 0", 2)
$d.Content.Find.Execute("372÷8=46, 4", $true, $false, $false, $false, $false, $true, 1, $false, "990÷8=123, 6", 2)
$d.Content.Find.Execute("330÷6=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "466÷7=66, 4", 2)
$d.Content.Find.Execute("921÷3=307, 0", $true, $false, $false, $false, $false, $true, 1, $false, "590÷2=295, 0", 2)
$d.Content.Find.Execute("362÷6=60, 2", $true, $false, $false, $false, $false, $true, 1, $false, "618÷5=123, 3", 2)
$d.Content.Find.Execute("752÷6=125, 2", $true, $false, $false, $false, $false, $true, 1, $false, "688÷5=137, 3", 2)
$d.Content.Find.Execute("303÷3=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "598÷4=149, 2", 2)
$d.Content.Find.Execute("184÷5=36, 4", $true, $false, $false, $false, $false, $true, 1, $false, "495÷7=70, 5", 2)
$d.Content.Find.Execute("694÷7=99, 1", $true, $false, $false, $false, $false, $true, 1, $false, "222÷5=44, 2", 2)
$d.Content.Find.Execute("164÷9=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "122÷5=24, 2", 2)
$d.Content.Find.Execute("727÷9=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "468÷3=156, 0", 2)
$d.Content.Find.Execute("457÷6=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "661÷2=330, 1", 2)
$d.Content.Find.Execute("430÷5=86, 0", $true, $false, $false, $false, $false, $true, 1, $false, "627÷4=156, 3", 2)
$d.Content.Find.Execute("381÷5=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "237÷7=33, 6", 2)
$d.Content.Find.Execute("768÷9=85, 3", $true, $false, $false, $false, $false, $true, 1, $false, "206÷7=29, 3", 2)
$d.Content.Find.Execute("756÷9=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "732÷4=183, 0", 2)
$d.Content.Find.Execute("608÷4=152, 0", $true, $false, $false, $false, $false, $true, 1, $false, "158÷4=39, 2", 2)
$d.Content.Find.Execute("478÷3=159, 1", $true, $false, $false, $false, $false, $true, 1, $false, "497÷6=82, 5", 2)
$d.Content.Find.Execute("219÷5=43, 4", $true, $false, $false, $false, $false, $true, 1, $false, "537÷2=268, 1", 2)
